$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 250; this shifts the existing rows 250-270
# down to 251-271 and carries formatting along (matches the diff, where
# every old row 250..270 reappears unchanged one row lower, and a brand
# new record shows up at row 250).
$ws.Rows(250).Insert()

# Populate the newly inserted row 250 with the new weekly record.
$ws.Range("A250").Value = 10
$ws.Range("B250").Value = "Vega Modelo de Temuco"
$ws.Range("C250").Value = "La Araucanía"
$ws.Range("D250").Value = 45013
$ws.Range("E250").Value = 9
$ws.Range("F250").Value = 100114007
$ws.Range("G250").Value = "Jengibre"
$ws.Range("H250").Value = "Sin especificar"
$ws.Range("I250").Value = "Primera"
$ws.Range("J250").Value = 50
$ws.Range("K250").Value = 25000
$ws.Range("L250").Value = 25000
$ws.Range("M250").Value = 25000
$ws.Range("N250").Value = "$/caja 13 kilos"
$ws.Range("O250").Value = "Perú"
$ws.Range("P250").Value = 1923
$ws.Range("Q250").Value = 13
$ws.Range("R250").Value = "Hortaliza"
